$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$biArr = New-Object 'object[,]' 24,8
$biArr[0,0] = 4.981230423876468
$biArr[0,1] = 1.533380357151145
$biArr[0,2] = 0.0765853248849524
$biArr[0,3] = 1.359268881961739
$biArr[0,4] = 2.948547961839751
$biArr[0,5] = 0.0007876981006996155
$biArr[0,6] = 0.01175518264854031
$biArr[0,7] = 0.006613395466705363
$biArr[1,0] = 4.320021657265556
$biArr[1,1] = 1.324722100298629
$biArr[1,2] = 0.07095763349168038
$biArr[1,3] = 1.172073800671214
$biArr[1,4] = 2.613109361582872
$biArr[1,5] = 0.0007946891518557333
$biArr[1,6] = 0.007212128027621834
$biArr[1,7] = 0.002948885422485059
$biArr[2,0] = 3.915698898342782
$biArr[2,1] = 1.198322027119559
$biArr[2,2] = 0.06745113563161453
$biArr[2,3] = 1.058569938732447
$biArr[2,4] = 2.409144672161048
$biArr[2,5] = 0.0007990927253246388
$biArr[2,6] = 0.004916224201339192
$biArr[2,7] = 0.00148838123403916
$biArr[3,0] = 3.747607318179689
$biArr[3,1] = 1.147871087881583
$biArr[3,2] = 0.06584242137925855
$biArr[3,3] = 1.012576996458506
$biArr[3,4] = 2.323316459685529
$biArr[3,5] = 0.0008009303667910995
$biArr[3,6] = 0.00408997314903714
$biArr[3,7] = 0.001131603897055466
$biArr[4,0] = 3.715280378918862
$biArr[4,1] = 1.140358903001413
$biArr[4,2] = 0.06537209896362839
$biArr[4,3] = 1.004904273159696
$biArr[4,4] = 2.305276147110703
$biArr[4,5] = 0.0008012544326233542
$biArr[4,6] = 0.003955060321159998
$biArr[4,7] = 0.00116115291958252
$biArr[5,0] = 3.901266069596943
$biArr[5,1] = 1.19995177266685
$biArr[5,2] = 0.06687462319382576
$biArr[5,3] = 1.05779948839654
$biArr[5,4] = 2.397540478082703
$biArr[5,5] = 0.0007991637887946879
$biArr[5,6] = 0.004891676594446448
$biArr[5,7] = 0.001675678934604186
$biArr[6,0] = 4.736434825295305
$biArr[6,1] = 1.464118339770437
$biArr[6,2] = 0.07391235557678044
$biArr[6,3] = 1.29416445724037
$biArr[6,4] = 2.818429149018527
$biArr[6,5] = 0.0007901460370847985
$biArr[6,6] = 0.01005516924093632
$biArr[6,7] = 0.005372860819191416
$biArr[7,0] = 6.415911772227162
$biArr[7,1] = 1.994740460635967
$biArr[7,2] = 0.0883769866692532
$biArr[7,3] = 1.772201791939992
$biArr[7,4] = 3.683463379825696
$biArr[7,5] = 0.0007732046513021398
$biArr[7,6] = 0.02468761372549633
$biArr[7,7] = 0.01955913258151032
$biArr[8,0] = 7.602905785113592
$biArr[8,1] = 2.376245661793007
$biArr[8,2] = 0.09544131048739501
$biArr[8,3] = 2.03044398709352
$biArr[8,4] = 4.270257071634973
$biArr[8,5] = 0.0007616360504728184
$biArr[8,6] = 0.0378572013089622
$biArr[8,7] = 0.03511095460282565
$biArr[9,0] = 7.588332565320059
$biArr[9,1] = 2.345199510257487
$biArr[9,2] = 0.0728433864938367
$biArr[9,3] = 1.330751556135809
$biArr[9,4] = 3.982694589184518
$biArr[9,5] = 0.0007601541995707229
$biArr[9,6] = 0.05280914036437068
$biArr[9,7] = 0.037753114894854
$biArr[10,0] = 7.345618155442139
$biArr[10,1] = 2.23660020803942
$biArr[10,2] = 0.05644071730884548
$biArr[10,3] = 0.8211603905815821
$biArr[10,4] = 3.643944556141406
$biArr[10,5] = 0.0007609435587364955
$biArr[10,6] = 0.08744102875231619
$biArr[10,7] = 0.03644100085519231
$biArr[11,0] = 6.897314969269019
$biArr[11,1] = 2.067507920854155
$biArr[11,2] = 0.04305388161407464
$biArr[11,3] = 0.4304368035855077
$biArr[11,4] = 3.232595640091887
$biArr[11,5] = 0.0007635703527378856
$biArr[11,6] = 0.1385865610561012
$biArr[11,7] = 0.03234772729667501
$biArr[12,0] = 6.49179278763927
$biArr[12,1] = 1.92469935826449
$biArr[12,2] = 0.03545943499494442
$biArr[12,3] = 0.2304414329774218
$biArr[12,4] = 2.914477541059597
$biArr[12,5] = 0.0007662003112132409
$biArr[12,6] = 0.1846997451284835
$biArr[12,7] = 0.02857043298932282
$biArr[13,0] = 6.345712463304778
$biArr[13,1] = 1.87803954795487
$biArr[13,2] = 0.03377085357345422
$biArr[13,3] = 0.1903441791731133
$biArr[13,4] = 2.817604012370936
$biArr[13,5] = 0.0007672860252298415
$biArr[13,6] = 0.1961007734092988
$biArr[13,7] = 0.02717402711722627
$biArr[14,0] = 5.936119174746807
$biArr[14,1] = 1.756462432841943
$biArr[14,2] = 0.03384028511147363
$biArr[14,3] = 0.1804204042892827
$biArr[14,4] = 2.657050463300919
$biArr[14,5] = 0.0007716125292273823
$biArr[14,6] = 0.1794527841903033
$biArr[14,7] = 0.02184367076104721
$biArr[15,0] = 5.833874973497018
$biArr[15,1] = 1.737849896047749
$biArr[15,2] = 0.03750419642386404
$biArr[15,3] = 0.2639330924608174
$biArr[15,4] = 2.69897714991464
$biArr[15,5] = 0.0007735824492763888
$biArr[15,6] = 0.1402455617029688
$biArr[15,7] = 0.01969550886613902
$biArr[16,0] = 5.994937761584595
$biArr[16,1] = 1.804951456726997
$biArr[16,2] = 0.04650595151861836
$biArr[16,3] = 0.5012350193891635
$biArr[16,4] = 2.933583967417889
$biArr[16,5] = 0.000773548116443686
$biArr[16,6] = 0.08838109454863741
$biArr[16,7] = 0.01962906721980584
$biArr[17,0] = 6.340717363587203
$biArr[17,1] = 1.94458517355838
$biArr[17,2] = 0.06104873384579435
$biArr[17,3] = 0.9453177939100215
$biArr[17,4] = 3.302873377397788
$biArr[17,5] = 0.0007717270255375275
$biArr[17,6] = 0.04693725659048198
$biArr[17,7] = 0.02186821644847559
$biArr[18,0] = 7.246826672728787
$biArr[18,1] = 2.281277834259754
$biArr[18,2] = 0.0917043484040363
$biArr[18,3] = 1.956744648804346
$biArr[18,4] = 4.078312994721784
$biArr[18,5] = 0.0007648015064704452
$biArr[18,6] = 0.03398296660081845
$biArr[18,7] = 0.03094044142827457
$biArr[19,0] = 8.271384052501162
$biArr[19,1] = 2.616745136655311
$biArr[19,2] = 0.1023510655173041
$biArr[19,3] = 2.323912087395797
$biArr[19,4] = 4.642254762457441
$biArr[19,5] = 0.0007553058648042135
$biArr[19,6] = 0.04725742275376454
$biArr[19,7] = 0.0460517376191687
$biArr[20,0] = 8.934333451542443
$biArr[20,1] = 2.825058429238027
$biArr[20,2] = 0.1084892713896934
$biArr[20,3] = 2.514524627111911
$biArr[20,4] = 4.999520173465328
$biArr[20,5] = 0.0007493132598017731
$biArr[20,6] = 0.05638144862008332
$biArr[20,7] = 0.05706598819542208
$biArr[21,0] = 8.595242952112244
$biArr[21,1] = 2.7106406466051
$biArr[21,2] = 0.1058939781486643
$biArr[21,3] = 2.412708661291333
$biArr[21,4] = 4.821158926806987
$biArr[21,5] = 0.0007524584084325714
$biArr[21,6] = 0.05146333044048212
$biArr[21,7] = 0.05092019961416483
$biArr[22,0] = 7.302328049739117
$biArr[22,1] = 2.290898539287753
$biArr[22,2] = 0.09503829538671482
$biArr[22,3] = 2.034688665809696
$biArr[22,4] = 4.138735064769037
$biArr[22,5] = 0.0007645285058988202
$biArr[22,6] = 0.03472078643161858
$biArr[22,7] = 0.03097258645892609
$biArr[23,0] = 5.935627322275309
$biArr[23,1] = 1.853165442785212
$biArr[23,2] = 0.08346568978059565
$biArr[23,3] = 1.640566109694802
$biArr[23,4] = 3.427107054559997
$biArr[23,5] = 0.0007777875869371744
$biArr[23,6] = 0.0201211202158289
$biArr[23,7] = 0.01505530797214849

$ws.Range("B2:I25").Value = $biArr

$pArr = New-Object 'object[,]' 24,1
$pArr[0,0] = 1.097255617216376
$pArr[1,0] = 1.094481175163864
$pArr[2,0] = 1.092854804039931
$pArr[3,0] = 1.09065326676032
$pArr[4,0] = 1.088405492363435
$pArr[5,0] = 1.087684391494015
$pArr[6,0] = 1.089487591704398
$pArr[7,0] = 1.103558314488197
$pArr[8,0] = 1.091024460534172
$pArr[9,0] = 0.9030669903682167
$pArr[10,0] = 0.7777018069335782
$pArr[11,0] = 0.6847524645093088
$pArr[12,0] = 0.636617747162866
$pArr[13,0] = 0.6281378259490218
$pArr[14,0] = 0.651435611291852
$pArr[15,0] = 0.6925813884899128
$pArr[16,0] = 0.7668408410710015
$pArr[17,0] = 0.8739235503081701
$pArr[18,0] = 1.077553176369335
$pArr[19,0] = 1.106713656422215
$pArr[20,0] = 1.121950448346524
$pArr[21,0] = 1.119745363144574
$pArr[22,0] = 1.103740619409947
$pArr[23,0] = 1.090272130385642

$ws.Range("P2:P25").Value = $pArr
